$word.UserName = "Tillmann Taape"
$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("Pour adjouster un canon faulcé", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($rng.End - 1, $rng.End)
$c0 = $d.Comments.Add($target, "line oneline two")
